$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G25").Value = 583
$ws.Range("G26").Value = 2
$ws.Range("G29").Value = 10
